# Auto-generated Excel COM-interop script applying the cryptos.xlsx data refresh
# (commit: "Updated cryptos list on Fri Aug 11 09:35:13 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text-safe cells (Coin name / Link / Volume%) -----------------------
# These never collide with Excel's automatic number detection (names, URLs,
# and the "  +/-X.XX%  " volume strings keep their padding spaces), so a
# plain .Value assignment keeps them stored as text, same as the source file.
$textUpdates = @(
    @('E3', '  -0.17%  '),
    @('E5', '  -0.74%  '),
    @('E6', '  +0.68%  '),
    @('E7', '  +0.01%  '),
    @('E8', '  +0.09%  '),
    @('E9', '  -0.60%  '),
    @('E10', '  +0.63%  '),
    @('E11', '  +0.20%  '),
    @('E12', '  -5.13%  '),
    @('E13', '  -0.15%  '),
    @('E14', '  -0.19%  '),
    @('E15', '  +2.78%  '),
    @('E16', '  -1.07%  '),
    @('E17', '  -1.35%  '),
    @('E18', '  -0.76%  '),
    @('E19', '  -2.57%  '),
    @('E20', '  -0.45%  '),
    @('E21', '  +0.04%  '),
    @('E22', '  -0.51%  '),
    @('E23', '  +0.05%  '),
    @('E24', '  +1.39%  '),
    @('E25', '  +0.55%  '),
    @('E26', '  -0.75%  '),
    @('E27', '  -0.36%  '),
    @('E28', '  -0.80%  '),
    @('B29', 'Toncoin'),
    @('C29', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'),
    @('E29', '  -0.48%  '),
    @('B30', 'Hedera'),
    @('C30', 'https://coinranking.com/coin/jad286TjB+hedera-hbar'),
    @('E30', '  -3.09%  '),
    @('E31', '  +0.88%  '),
    @('E32', '  -0.61%  '),
    @('E33', '  -2.73%  '),
    @('E34', '  -1.26%  '),
    @('E35', '  -1.09%  '),
    @('E36', '  +0.08%  '),
    @('E37', '  +1.28%  '),
    @('E38', '  +1.54%  '),
    @('E39', '  -0.23%  '),
    @('E40', '  +0.47%  '),
    @('E41', '  +1.13%  '),
    @('E42', '  +0.05%  '),
    @('E43', '  -6.50%  '),
    @('E44', '  -0.63%  '),
    @('E45', '  -0.88%  '),
    @('E46', '  -3.59%  '),
    @('B47', 'BabyDogeCoin'),
    @('C47', 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'),
    @('E47', '  -0.06%  '),
    @('B48', 'TheSandbox'),
    @('C48', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'),
    @('E48', '  -0.18%  '),
    @('B49', 'EnergySwap'),
    @('C49', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'),
    @('E49', '  -0.84%  '),
    @('E50', '  -1.86%  '),
    @('E51', '  +1.03%  ')
)
foreach ($pair in $textUpdates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# --- Price cells (column D) ----------------------------------------------
# Many of the new prices are plain decimal numbers ("240.07", "1.000", ...).
# Assigning those directly would make Excel auto-convert the cell to a number
# (dropping significant trailing zeros / changing the stored type), so each
# price cell is briefly switched to Text format, written, then restored to
# the workbook's default "Normal" style so the saved file keeps the original
# (un-styled) General-format text cells.
$priceUpdates = @(
    @('D2', '29.359.45'),
    @('D3', '1.846.02'),
    @('D5', '240.07'),
    @('D6', '0.6308'),
    @('D8', '0.07535'),
    @('D9', '0.2957'),
    @('D10', '24.44'),
    @('D11', '0.07717'),
    @('D12', '1.846.43'),
    @('D13', '4.990'),
    @('D14', '0.6832'),
    @('D15', '0.000009997'),
    @('D16', '82.78'),
    @('D17', '6.130'),
    @('D18', '29.403.54'),
    @('D19', '227.45'),
    @('D20', '12.42'),
    @('D21', '1.000'),
    @('D22', '7.541'),
    @('D24', '157.28'),
    @('D25', '0.1397'),
    @('D28', '1.466'),
    @('D29', '1.256'),
    @('D30', '0.05684'),
    @('D31', '4.126'),
    @('D32', '4.014'),
    @('D33', '1.844'),
    @('D34', '1.154'),
    @('D35', '0.7142'),
    @('D36', '2.590'),
    @('D37', '1.255.09'),
    @('D39', '2.786'),
    @('D40', '0.9124'),
    @('D41', '6.205'),
    @('D43', '2.001.82'),
    @('D44', '101.22'),
    @('D45', '66.45'),
    @('D46', '7.041'),
    @('D47', '0.00000000117'),
    @('D48', '0.4021'),
    @('D49', '9.088'),
    @('D50', '1.690'),
    @('D51', '0.1125')
)
foreach ($pair in $priceUpdates) {
    $cell = $ws.Range($pair[0])
    $cell.NumberFormat = "@"
    $cell.Value = $pair[1]
    $cell.Style = "Normal"
}

Write-Output "Applied $($textUpdates.Count) text updates and $($priceUpdates.Count) price updates."
